$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (PAN), shifting PAN..Send Confirmation Email right by one.
$ws.Columns("F").Insert()

# Header for the new "Phone" column.
$ws.Range("F1").Value = "Phone"

# Phone values for the two investor rows.
$ws.Range("F2").Value = 999999999
$ws.Range("F3").Value = 111111111

# Match the hyperlink-like style (underline + hyperlink theme color) used on
# the phone cells in the target workbook.
$ws.Range("F2:F3").Style = "Hyperlink"

# Column width for the new column: narrower, not best-fit like its
# neighbours -- matches the email column's (manually-set) width.
$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth

# Restore the selection Excel leaves behind after this kind of edit.
$ws.Range("F4").Select()
